$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 corresponds to year 2004. Update values to match corrected
# "silber 93-04" series (subtle differences fix).
$ws.Range("C13").Value = 230361
$ws.Range("D13").Value = 52429.3
$ws.Range("E13").Value = 1.539949830059086
$ws.Range("J13").Value = 39436.7
$ws.Range("K13").Value = 1695099
$ws.Range("L13").Value = 0.6480785059764926
$ws.Range("M13").Value = 1.13974394575808
$ws.Range("N13").Value = 0.5686176341524102
$ws.Range("O13").Value = 0.5707149543900899
$ws.Range("P13").Value = 0.1196394627777258
